$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13226
$ws1.Range("F5").Value = 10
$ws1.Range("F6").Value = 110
$ws1.Range("F11").Value = 13167
$ws1.Range("F12").Value = 319
$ws1.Range("F13").Value = 573
$ws1.Range("F14").Value = 8830
$ws1.Range("F15").Value = 7902
$ws1.Range("F16").Value = 227
$ws1.Range("F20").Value = 2
$ws1.Range("F21").Value = 9
$ws1.Range("F22").Value = 1000
$ws1.Range("F28").Value = 350

# Sheet "演出" (sheet2) updates to column F
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 25

# Sheet "全部类型" (sheet4) updates to column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13226
$ws4.Range("F6").Value = 10
$ws4.Range("F7").Value = 110
$ws4.Range("F12").Value = 13167
$ws4.Range("F13").Value = 319
$ws4.Range("F14").Value = 573
$ws4.Range("F15").Value = 8830
$ws4.Range("F16").Value = 7902
$ws4.Range("F17").Value = 227
$ws4.Range("F21").Value = 2
$ws4.Range("F22").Value = 9
$ws4.Range("F23").Value = 1000
$ws4.Range("F26").Value = 25
$ws4.Range("F31").Value = 350
